$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Cell values -----------------------------------------------------
# Row 1
$ws.Range("A1").Value2 = "Username"
$ws.Range("B1").Value2 = "Label"
$ws.Range("C1").Value2 = 1
$ws.Range("D1").Value2 = "a"

# Row 2
$ws.Range("A2").Value2 = "Username"
$ws.Range("B2").Value2 = "Text"
$ws.Range("C2").Value2 = 2
$ws.Range("D2").Value2 = "[a-zA-Z0-9]*"

# Row 3
$ws.Range("A3").Value2 = "Email"
$ws.Range("B3").Value2 = "Label"
$ws.Range("C3").Value2 = 3
$ws.Range("D3").Value2 = "a"

# Row 4
$ws.Range("A4").Value2 = "Email"
$ws.Range("B4").Value2 = "Email"
$ws.Range("C4").Value2 = 4
$ws.Range("D4").Value2 = "[a-zA-Z0-9]*@[a-zA-Z0-9]*"

# Row 5
$ws.Range("A5").Value2 = "Confirm Email"
$ws.Range("B5").Value2 = "Email"
$ws.Range("C5").Value2 = 5
$ws.Range("D5").Value2 = "[a-zA-Z0-9]*@[a-zA-Z0-9]*"
$ws.Range("E5").Value2 = "Eq(4)"

# Row 6
$ws.Range("A6").Value2 = "Password"
$ws.Range("B6").Value2 = "Label"
$ws.Range("C6").Value2 = 6
$ws.Range("D6").Value2 = "a"

# Row 7
$ws.Range("A7").Value2 = "Password"
$ws.Range("B7").Value2 = "Password"
$ws.Range("C7").Value2 = 7
$ws.Range("D7").Value2 = "[a-zA-Z0-9]*"

# ---- Formatting --------------------------------------------------------
# Every populated cell in the used range carries the same (non-default)
# cell style after the edit - apply it row by row so we don't spill
# formatting onto untouched columns (e.g. E1:E4, E6:E7).
$ws.Range("A1:D1").Style = "Normal"
$ws.Range("A2:D2").Style = "Normal"
$ws.Range("A3:D3").Style = "Normal"
$ws.Range("A4:D4").Style = "Normal"
$ws.Range("A5:E5").Style = "Normal"
$ws.Range("A6:D6").Style = "Normal"
$ws.Range("A7:D7").Style = "Normal"

# ---- Column widths -------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.833333333333332
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(4).ColumnWidth = 40.833333333333336

# ---- Selection -------------------------------------------------------
$ws.Range("F11").Select() | Out-Null
